# Delete the "St Clements Hospital" row (row 4) from Sheet1.
# This removes the row entirely, shifting all subsequent rows up by one,
# and removes the now-unused shared strings for that row's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(4).Delete()

# Mirror the selection state Excel leaves behind after deleting an entire row:
# the row that shifted into the deleted row's place becomes selected in full.
$ws.Rows.Item(4).Select() | Out-Null
